# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F3").Value = 327
$sheet1.Range("F4").Value = 418
$sheet1.Range("F7").Value = 2161
$sheet1.Range("F11").Value = 4831
$sheet1.Range("F17").Value = 170
$sheet1.Range("F20").Value = 116
$sheet1.Range("F21").Value = 3761
$sheet1.Range("F22").Value = 692
$sheet1.Range("F23").Value = 627
$sheet1.Range("F28").Value = 18
$sheet1.Range("F31").Value = 571
$sheet1.Range("F33").Value = 23
$sheet1.Range("F34").Value = 881
$sheet1.Range("F35").Value = 2393

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F3").Value = 327
$sheet4.Range("F4").Value = 418
$sheet4.Range("F7").Value = 2161
$sheet4.Range("F11").Value = 4831
$sheet4.Range("F17").Value = 170
$sheet4.Range("F20").Value = 116
$sheet4.Range("F21").Value = 3761
$sheet4.Range("F22").Value = 692
$sheet4.Range("F23").Value = 627
$sheet4.Range("F28").Value = 18
$sheet4.Range("F31").Value = 571
$sheet4.Range("F34").Value = 23
$sheet4.Range("F35").Value = 881
$sheet4.Range("F36").Value = 2393
